$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Add the two new header columns for per-line pricing.
$ws.Range("E1").Value = "Unit Price"
$ws.Range("F1").Value = "Total Price"

# Match the author's final selection state (cell F2 selected).
$ws.Activate()
[void]$ws.Range("F2").Select()
